$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H15").Value = 1148.6724
$ws.Range("I15").Value = 1148.6724
$ws.Range("K15").Value = 3446.0172
$ws.Range("M15").Value = -3277.0172

$ws.Range("H33").Value = 147.25
$ws.Range("I33").Value = 111.333336
$ws.Range("K33").Value = 111.333336
$ws.Range("M33").Value = 117.666664

$ws.Range("H62").Value = 9269.25
$ws.Range("I62").Value = 8105.3125
$ws.Range("J62").Value = 13925
$ws.Range("K62").Value = 8105.3125
$ws.Range("L62").Value = 13925
$ws.Range("M62").Value = -7481.3125
$ws.Range("N62").Value = -15173

$ws.Range("H65").Value = 9269.25
$ws.Range("I65").Value = 8105.3125
$ws.Range("J65").Value = 13925
$ws.Range("K65").Value = 40526.5625
$ws.Range("L65").Value = 69625
$ws.Range("M65").Value = -37406.5625
$ws.Range("N65").Value = -75865

$ws.Range("H95").Value = 10000
$ws.Range("J95").Value = 10000
$ws.Range("L95").Value = 10000
$ws.Range("N95").Value = -15492

$ws.Range("H98").Value = 96155420
$ws.Range("I98").Value = 96155420
$ws.Range("K98").Value = 96155420
$ws.Range("M98").Value = -96153922

$ws.Range("H113").Value = 2824.4736
$ws.Range("I113").Value = 2836.5
$ws.Range("J113").Value = 2811.111
$ws.Range("K113").Value = 2836.5
$ws.Range("L113").Value = 2811.111
$ws.Range("M113").Value = 417.5
$ws.Range("N113").Value = -9319.111000000001

$ws.Range("H122").Value = 96155420
$ws.Range("I122").Value = 96155420
$ws.Range("K122").Value = 288466260
$ws.Range("M122").Value = -288463810

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H5").Value = 158.88889
$ws.Range("I5").Value = 155.2
$ws.Range("J5").Value = 163.5
$ws.Range("K5").Value = 155.2
$ws.Range("L5").Value = 163.5
$ws.Range("M5").Value = -43.19999999999999
$ws.Range("N5").Value = -387.5

$ws.Range("H32").Value = 7579257
$ws.Range("I32").Value = 3072.5366
$ws.Range("J32").Value = 111120450
$ws.Range("K32").Value = 3072.5366
$ws.Range("L32").Value = 111120450
$ws.Range("M32").Value = -2785.5366
$ws.Range("N32").Value = -111121024

$ws.Range("H45").Value = 72854.42999999999
$ws.Range("I45").Value = 126267.25
$ws.Range("J45").Value = 1637.3334
$ws.Range("K45").Value = 126267.25
$ws.Range("L45").Value = 1637.3334
$ws.Range("M45").Value = -125890.25
$ws.Range("N45").Value = -2391.3334

$ws.Range("H132").Value = 2263992
$ws.Range("I132").Value = 961
$ws.Range("J132").Value = 5884842
$ws.Range("K132").Value = 2883
$ws.Range("L132").Value = 17654526
$ws.Range("M132").Value = -353
$ws.Range("N132").Value = -17659586

$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H4").Value = 158.88889
$ws.Range("I4").Value = 155.2
$ws.Range("J4").Value = 163.5
$ws.Range("K4").Value = 155.2
$ws.Range("L4").Value = 163.5
$ws.Range("M4").Value = -40.19999999999999
$ws.Range("N4").Value = -393.5

$ws.Range("H76").Value = 21771.334
$ws.Range("I76").Value = 0
$ws.Range("J76").Value = 21771.334
$ws.Range("K76").Value = 0
$ws.Range("L76").Value = 21771.334
$ws.Range("M76").ClearContents()
$ws.Range("N76").Value = -22401.334

$ws.Range("H79").Value = 21771.334
$ws.Range("I79").Value = 0
$ws.Range("J79").Value = 21771.334
$ws.Range("K79").Value = 0
$ws.Range("L79").Value = 21771.334
$ws.Range("M79").ClearContents()
$ws.Range("N79").Value = -23955.334

$ws.Range("H82").Value = 7730
$ws.Range("I82").Value = 4936.7144
$ws.Range("J82").Value = 27283
$ws.Range("K82").Value = 4936.7144
$ws.Range("L82").Value = 27283
$ws.Range("M82").Value = -4553.7144
$ws.Range("N82").Value = -28049

$ws.Range("H85").Value = 7730
$ws.Range("I85").Value = 4936.7144
$ws.Range("J85").Value = 27283
$ws.Range("K85").Value = 4936.7144
$ws.Range("L85").Value = 27283
$ws.Range("M85").Value = -3610.7144
$ws.Range("N85").Value = -29935

$ws.Range("H120").Value = 0
$ws.Range("J120").Value = 0
$ws.Range("L120").Value = 0
$ws.Range("N120").ClearContents()

$ws.Range("H134").Value = 7416615
$ws.Range("I134").Value = 2868.4443
$ws.Range("J134").Value = 18537236
$ws.Range("K134").Value = 8605.332900000001
$ws.Range("L134").Value = 55611708
$ws.Range("M134").Value = -6070.332900000001
$ws.Range("N134").Value = -55616778

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H6").Value = 1066.3334
$ws.Range("I6").Value = 199
$ws.Range("K6").Value = 199
$ws.Range("M6").Value = -86

$ws.Range("H7").Value = 161.5
$ws.Range("I7").Value = 161.7
$ws.Range("J7").Value = 161.3
$ws.Range("K7").Value = 161.7
$ws.Range("L7").Value = 161.3
$ws.Range("M7").Value = -48.69999999999999
$ws.Range("N7").Value = -387.3

$ws.Range("H16").Value = 1437.5714
$ws.Range("I16").Value = 918.75
$ws.Range("J16").Value = 2129.3333
$ws.Range("K16").Value = 918.75
$ws.Range("L16").Value = 2129.3333
$ws.Range("M16").Value = -631.75
$ws.Range("N16").Value = -2703.3333

$ws.Range("H18").Value = 44000
$ws.Range("J18").Value = 44000
$ws.Range("L18").Value = 44000
$ws.Range("N18").Value = -44460

$ws.Range("H25").Value = 7500
$ws.Range("J25").Value = 10000
$ws.Range("L25").Value = 10000
$ws.Range("N25").Value = -10348

$ws.Range("H55").Value = 0
$ws.Range("I55").Value = 0
$ws.Range("K55").Value = 0
$ws.Range("M55").ClearContents()

$ws.Range("H99").Value = 125002130
$ws.Range("I99").Value = 250002320
$ws.Range("J99").Value = 1925
$ws.Range("K99").Value = 250002320
$ws.Range("L99").Value = 1925
$ws.Range("M99").Value = -250000822
$ws.Range("N99").Value = -4921

$ws.Range("H113").Value = 1437.5714
$ws.Range("I113").Value = 918.75
$ws.Range("J113").Value = 2129.3333
$ws.Range("K113").Value = 918.75
$ws.Range("L113").Value = 2129.3333
$ws.Range("M113").Value = 1251.25
$ws.Range("N113").Value = -6469.3333

$ws.Range("H117").Value = 0
$ws.Range("J117").Value = 0
$ws.Range("L117").Value = 0
$ws.Range("N117").ClearContents()

$ws.Range("H118").Value = 44850
$ws.Range("J118").Value = 44850
$ws.Range("L118").Value = 44850
$ws.Range("N118").Value = -48164

$ws.Range("H126").Value = 125002130
$ws.Range("I126").Value = 250002320
$ws.Range("J126").Value = 1925
$ws.Range("K126").Value = 750006960
$ws.Range("L126").Value = 5775
$ws.Range("M126").Value = -750004490
$ws.Range("N126").Value = -10715

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H131").Value = 838.1799999999999
$ws.Range("J131").Value = 864.56384
$ws.Range("L131").Value = 2593.69152
$ws.Range("N131").Value = -12673.69152

$ws.Range("H140").Value = 10001130
$ws.Range("I140").Value = 10870650
$ws.Range("K140").Value = 32611950
$ws.Range("M140").Value = -32606770

$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H80").Value = 3596.25
$ws.Range("I80").Value = 4108.75
$ws.Range("J80").Value = 2827.5
$ws.Range("K80").Value = 4108.75
$ws.Range("L80").Value = 2827.5
$ws.Range("M80").Value = -3110.75
$ws.Range("N80").Value = -4823.5

$ws.Range("H83").Value = 3596.25
$ws.Range("I83").Value = 4108.75
$ws.Range("J83").Value = 2827.5
$ws.Range("K83").Value = 20543.75
$ws.Range("L83").Value = 14137.5
$ws.Range("M83").Value = -15551.75
$ws.Range("N83").Value = -24121.5

$ws.Range("H93").Value = 22000
$ws.Range("J93").Value = 22000
$ws.Range("L93").Value = 22000
$ws.Range("N93").Value = -25744

$ws.Range("H107").Value = 385.25
$ws.Range("I107").Value = 231.66667
$ws.Range("J107").Value = 846
$ws.Range("K107").Value = 231.66667
$ws.Range("L107").Value = 846
$ws.Range("M107").Value = 1688.33333
$ws.Range("N107").Value = -4686

$ws.Range("H113").Value = 2941
$ws.Range("I113").Value = 3000
$ws.Range("J113").Value = 2911.5
$ws.Range("K113").Value = 3000
$ws.Range("L113").Value = 2911.5
$ws.Range("M113").Value = -830
$ws.Range("N113").Value = -7251.5

$ws.Range("H132").Value = 10763
$ws.Range("I132").Value = 2136
$ws.Range("J132").Value = 20827.834
$ws.Range("K132").Value = 6408
$ws.Range("L132").Value = 62483.50199999999
$ws.Range("M132").Value = -3878
$ws.Range("N132").Value = -67543.50199999999

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H61").Value = 2050.75
$ws.Range("I61").Value = 2050.75
$ws.Range("K61").Value = 2050.75
$ws.Range("M61").Value = -1848.75

$ws.Range("H113").Value = 2050.75
$ws.Range("I113").Value = 2050.75
$ws.Range("K113").Value = 2050.75
$ws.Range("M113").Value = 119.25

$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H92").Value = 53383.332
$ws.Range("J92").Value = 53383.332
$ws.Range("L92").Value = 53383.332
$ws.Range("N92").Value = -58375.332
